$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "object" row appended to the RGB table: AusgangRechts = (255, 20, 0)
$ws.Range("A25").Value = "AusgangRechts"
$ws.Range("B25").Value = 255
$ws.Range("C25").Value = 20
$ws.Range("D25").Value = 0

# Keep the sheet view/selection in sync with the newly added row, like
# Excel does after typing into the row just past the previous selection.
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D26").Select()
